# Update "paises.xlsx" worksheet (Pais) with refreshed COVID-19 country data
# and the new "Datos actualizados" timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "last updated" timestamp banner (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 5 de Abril de 2020 a las 00:52"

# --- Row 4: Estados Unidos ---
$ws.Range("B4").Value = 309405
$ws.Range("C4").Value = 32244
$ws.Range("D4").Value = 14741
$ws.Range("E4").Value = 286234
$ws.Range("F4").Value = 8206
$ws.Range("G4").Value = 1026
$ws.Range("H4").Value = 8430

# --- Row 19 ---
$ws.Range("B19").Value = 10323
$ws.Range("C19").Value = 1129
$ws.Range("E19").Value = 9757
$ws.Range("G19").Value = 76
$ws.Range("H19").Value = 439

# --- Row 27 ---
$ws.Range("B27").Value = 4472
$ws.Range("C27").Value = 282
$ws.Range("E27").Value = 4335

# --- Rows 51-53: Argentina moves above Islandia and Colombia ---
$ws.Range("A51").Value = "Argentina"
$ws.Range("B51").Value = 1451
$ws.Range("C51").Value = 98
$ws.Range("D51").Value = 279
$ws.Range("E51").Value = 1129
$ws.Range("F51").Value = 0
$ws.Range("G51").Value = 1
$ws.Range("H51").Value = 43

$ws.Range("A52").Value = "Islandia"
$ws.Range("B52").Value = 1417
$ws.Range("C52").Value = 53
$ws.Range("D52").Value = 396
$ws.Range("E52").Value = 1017
$ws.Range("F52").Value = 12
$ws.Range("G52").Value = 0
$ws.Range("H52").Value = 4

$ws.Range("A53").Value = "Colombia"
$ws.Range("B53").Value = 1406
$ws.Range("C53").Value = 139
$ws.Range("D53").Value = 85
$ws.Range("E53").Value = 1289
$ws.Range("F53").Value = 50
$ws.Range("G53").Value = 7
$ws.Range("H53").Value = 32

# --- Rows 133-134: Guayana Francesa moves above Guatemala ---
$ws.Range("A133").Value = "Guayana Francesa"
$ws.Range("B133").Value = 61
$ws.Range("C133").Value = 4
$ws.Range("D133").Value = 22
$ws.Range("E133").Value = 39
$ws.Range("G133").Value = 0
$ws.Range("H133").Value = 0

$ws.Range("A134").Value = "Guatemala"
$ws.Range("C134").Value = 7
$ws.Range("D134").Value = 15
$ws.Range("E134").Value = 40
$ws.Range("G134").Value = 1
$ws.Range("H134").Value = 2

# --- Row 140 ---
$ws.Range("B140").Value = 44
$ws.Range("C140").Value = 2
$ws.Range("E140").Value = 34

# --- Row 151 ---
$ws.Range("B151").Value = 28
$ws.Range("C151").Value = 4
$ws.Range("D151").Value = 0
$ws.Range("E151").Value = 24
$ws.Range("G151").Value = 1
$ws.Range("H151").Value = 4

# --- Rows 162-163: Libia moves above Islas Virgenes de los Estados Unidos ---
$ws.Range("A162").Value = "Libia"
$ws.Range("B162").Value = 18
$ws.Range("C162").Value = 1
$ws.Range("H162").Value = 1

$ws.Range("A163").Value = "Islas Virgenes de los Estados Unidos"
$ws.Range("E163").Value = 17
$ws.Range("H163").Value = 0
